$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new, empty paragraph right before the "A Nice description..."
#    paragraph (the second paragraph in the body).
# ------------------------------------------------------------------
$descParaBefore = $d.Paragraphs(2)
$descParaBefore.Range.InsertParagraphBefore()

# InsertParagraphBefore() leaves a stray empty run behind on the newly
# created paragraph (<w:p><w:r></w:r></w:p>); replace its (empty) range
# with clean OOXML so it serializes as a bare <w:p/>.
$newEmptyPara = $d.Paragraphs(2)
$emptyRange = $newEmptyPara.Range
$emptyParaXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$emptyRange.InsertXML($emptyParaXml)

# ------------------------------------------------------------------
# 2. Re-split the description paragraph's text into four runs:
#      "A "  /  "n"  /  "ice description...pictures"  /  "."
#    (i.e. "A Nice description..." -> "A nice description..."),
#    keeping the paragraph's own formatting/identity untouched.
# ------------------------------------------------------------------
$descPara = $d.Paragraphs(3)
$descRange = $descPara.Range
# Paragraph.Range includes the trailing paragraph mark; trim it off so
# InsertXML only replaces the run content, not the paragraph mark itself.
$textRange = $d.Range($descRange.Start, $descRange.End - 1)

$splitXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">A </w:t></w:r><w:r><w:t>n</w:t></w:r><w:r><w:t>ice description of the college. This can be multiple paragraphs long and include pictures</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$textRange.InsertXML($splitXml)
